$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Merge row 3's data into row 2 (Species + Unknown Sex Count),
# then remove row 3 entirely.
$ws.Range("F2").Value = "Fisher"
$ws.Range("J2").Value = 1

$ws.Rows.Item(3).Delete()
